$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C text updates (coin name / link swap for rows 43-44) ---
$ws.Range('B43').Value = 'Frax'
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

# --- Column D price updates. These cells store plain text (e.g. "1.001",
# "22.388.07") in the workbook, but Excel's COM Value setter auto-parses
# number-looking strings into real numbers. To keep them as literal text
# (matching the original inline-string cell type) we briefly force a text
# number format, assign the value, then restore the "Normal" style so no
# extra formatting is left behind.
$dCells = @{
    'D2' = '22.388.07'
    'D3' = '1.560.29'
    'D6' = '285.82'
    'D9' = '0.3332'
    'D11' = '0.07382'
    'D15' = '6.863'
    'D16' = '1.560.81'
    'D18' = '88.63'
    'D19' = '0.06710'
    'D21' = '6.307'
    'D22' = '15.99'
    'D23' = '11.89'
    'D24' = '22.386.00'
    'D25' = '2.384'
    'D26' = '2.535'
    'D27' = '149.13'
    'D29' = '5.003'
    'D30' = '122.75'
    'D31' = '1.734.69'
    'D33' = '6.086'
    'D34' = '1.985'
    'D36' = '0.08220'
    'D39' = '0.2204'
    'D40' = '0.06333'
    'D41' = '5.293'
    'D42' = '11.07'
    'D43' = '1.001'
    'D44' = '0.6036'
    'D45' = '13.63'
    'D46' = '3.756'
    'D47' = '0.5709'
    'D48' = '124.09'
    'D50' = '1.207'
    'D51' = '0.07211'
}
foreach ($addr in $dCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dCells[$addr]
    $cell.Style = "Normal"
}

# --- Column E volume-percentage updates (plain text, padded with spaces) ---
$eCells = @{
    'E2' = '  -0.34%  '
    'E3' = '  -0.83%  '
    'E4' = '  -0.10%  '
    'E6' = '  -1.87%  '
    'E7' = '  -2.81%  '
    'E8' = '  -2.43%  '
    'E9' = '  -2.15%  '
    'E10' = '  -2.01%  '
    'E11' = '  -2.55%  '
    'E12' = '  -0.07%  '
    'E13' = '  -3.09%  '
    'E14' = '  -1.41%  '
    'E15' = '  -1.27%  '
    'E16' = '  -0.69%  '
    'E17' = '  -2.30%  '
    'E18' = '  -2.70%  '
    'E19' = '  -0.36%  '
    'E20' = '  -0.11%  '
    'E21' = '  +0.32%  '
    'E22' = '  -2.76%  '
    'E23' = '  -2.63%  '
    'E24' = '  -0.36%  '
    'E25' = '  +2.43%  '
    'E26' = '  -2.41%  '
    'E27' = '  +0.46%  '
    'E28' = '  -3.98%  '
    'E29' = '  -0.02%  '
    'E30' = '  -2.58%  '
    'E31' = '  -0.70%  '
    'E32' = '  +0.42%  '
    'E33' = '  -1.14%  '
    'E34' = '  +0.04%  '
    'E35' = '  -3.22%  '
    'E36' = '  -2.68%  '
    'E37' = '  -4.04%  '
    'E38' = '  -6.70%  '
    'E39' = '  -4.14%  '
    'E40' = '  -3.30%  '
    'E41' = '  -3.84%  '
    'E42' = '  -2.89%  '
    'E43' = '  -0.08%  '
    'E44' = '  -4.02%  '
    'E45' = '  -2.41%  '
    'E46' = '  -1.52%  '
    'E47' = '  -2.80%  '
    'E48' = '  -4.48%  '
    'E49' = '  -4.75%  '
    'E50' = '  -2.01%  '
    'E51' = '  -1.62%  '
}
foreach ($addr in $eCells.Keys) {
    $ws.Range($addr).Value = $eCells[$addr]
}
